$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new data rows for 2025-10-11 (serial 45941)
$ws.Range("A22").Value = 45941
$ws.Range("B22").Value = "四方坪站"
$ws.Range("C22").Value = 10508.12
$ws.Range("D22").Value = 8574.84
$ws.Range("E22").Value = 3652.19
$ws.Range("F22").Value = 446

$ws.Range("A23").Value = 45941
$ws.Range("B23").Value = "高岭站"
$ws.Range("C23").Value = 5191.68
$ws.Range("D23").Value = 4116
$ws.Range("E23").Value = 1371.67
$ws.Range("F23").Value = 189

# Update the visible selection to H20 (also clears the stale topLeftCell scroll position)
$ws.Range("H20").Select() | Out-Null
